# employees.xlsx automation:
# - Sheet "Sheet" holds the list of employee names (A1, A2, ...)
# - For each employee, create a dedicated sheet named after them with a
#   "Hello <name>!" greeting in A1.

$wb = $excel.ActiveWorkbook

# Sheet 1 ("Sheet") now carries the employee roster, replacing the old
# "Kitty" sheet's hello/world placeholder content.
$names = @("Maeva", "Gertude")

$rosterSheet = $wb.Worksheets.Item(1)
for ($i = 0; $i -lt $names.Count; $i++) {
    $rosterSheet.Cells.Item($i + 1, 1).Value = $names[$i]
}

# Rename the old placeholder "Kitty" sheet to the first employee and give
# it its greeting (clearing the old hello/world placeholder content first).
$firstSheet = $wb.Worksheets.Item(2)
$firstSheet.Name = $names[0]
$firstSheet.Cells.Clear()
$firstSheet.Range("A1").Value = "Hello " + $names[0] + "!"

# Append a new sheet per remaining employee, each with its own greeting.
for ($i = 1; $i -lt $names.Count; $i++) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $names[$i]
    $newSheet.Range("A1").Value = "Hello " + $names[$i] + "!"
}
